$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.508.30"
$ws.Range("E2").Value = "  -3.18%  "
$ws.Range("D3").Value = "1.659.81"
$ws.Range("E3").Value = "  -3.96%  "
$ws.Range("E4").Value = "  +0.04%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "214.50"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  -1.98%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.511"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  -2.42%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.999"
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = "  +0.04%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "24.45"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "  +1.90%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.263"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  -1.67%  "
$ws.Range("E10").Value = "  -2.53%  "
$ws.Range("E11").Value = "  -1.77%  "
$ws.Range("D12").Value = "1.894.91"
$ws.Range("E12").Value = "  -3.84%  "
$ws.Range("D13").Value = "1.657.22"
$ws.Range("E13").Value = "  -4.05%  "
$ws.Range("E15").Value = "  +0.19%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "65.89"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "  -2.52%  "
$ws.Range("D17").Value = "27.518.47"
$ws.Range("E17").Value = "  -2.86%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "240.86"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  -2.22%  "
$ws.Range("E19").Value = "  -2.79%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "7.63"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  -3.41%  "
$ws.Range("E21").Value = "  +0.09%  "
$ws.Range("E22").Value = "  -3.51%  "
$ws.Range("E23").Value = "  -2.49%  "
$ws.Range("E24").Value = "  -1.44%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "146.02"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  -2.21%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "7.23"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  -2.87%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "16.28"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  -1.97%  "
$ws.Range("E28").Value = "  -0.08%  "
$ws.Range("E29").Value = "  -2.26%  "
$ws.Range("E30").Value = "  -3.08%  "
$ws.Range("E31").Value = "  -1.09%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.33"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  -2.84%  "
$ws.Range("D33").Value = "1.453.86"
$ws.Range("E33").Value = "  -2.19%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "3.10"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  -4.74%  "
$ws.Range("E35").Value = "  -4.41%  "
$ws.Range("E36").Value = "  -1.30%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.923"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  -5.72%  "
$ws.Range("E38").Value = "  -4.61%  "
$ws.Range("E39").Value = "  -2.97%  "
$ws.Range("E40").Value = "  -0.82%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.00"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  -0.04%  "
$ws.Range("E42").Value = "  -4.26%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "5.46"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  -3.32%  "
$ws.Range("E44").Value = "  -2.72%  "
$ws.Range("B45").Value = "TrustWalletToken"
$ws.Range("C45").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.790"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  -2.02%  "
$ws.Range("B46").Value = "RocketPoolETH"
$ws.Range("C46").Value = "https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"
$ws.Range("D46").Value = "1.802.41"
$ws.Range("E46").Value = "  -3.91%  "
$ws.Range("E47").Value = "  -0.35%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "88.61"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  -1.91%  "
$ws.Range("E49").Value = "  -5.87%  "
$ws.Range("E50").Value = "  -1.43%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "7.84"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  -3.68%  "
